$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='60.337.49'; E='  +5.86%  ' },
  @{ Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='2.654.51'; E='  +10.15%  ' },
  @{ Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='''0.999'; E='  -0.21%  ' },
  @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='''507.82'; E='  +4.09%  ' },
  @{ Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='''157.26'; E='  +2.45%  ' },
  @{ Row=7; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='''0.995'; E='  -0.23%  ' },
  @{ Row=8; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='''0.591'; E='  -4.00%  ' },
  @{ Row=9; B='LidoStakedEther'; C='https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'; D='2.650.22'; E='  +9.28%  ' },
  @{ Row=10; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='''6.49'; E='  +3.72%  ' },
  @{ Row=11; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='''0.105'; E='  +5.34%  ' },
  @{ Row=12; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='''0.344'; E='  +3.55%  ' },
  @{ Row=13; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='''0.127'; E='  +0.88%  ' },
  @{ Row=14; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='3.080.31'; E='  +8.67%  ' },
  @{ Row=15; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='60.473.51'; E='  +6.10%  ' },
  @{ Row=16; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='''21.67'; E='  +5.64%  ' },
  @{ Row=17; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='''0.0000140'; E='  +6.02%  ' },
  @{ Row=18; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='2.646.29'; E='  +8.97%  ' },
  @{ Row=19; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='''4.79'; E='  +2.61%  ' },
  @{ Row=20; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='''345.86'; E='  +7.09%  ' },
  @{ Row=21; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='''10.50'; E='  +5.27%  ' },
  @{ Row=22; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='''6.19'; E='  +4.66%  ' },
  @{ Row=23; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='''0.998'; E='  +0.24%  ' },
  @{ Row=24; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='''59.97'; E='  +4.11%  ' },
  @{ Row=25; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='''0.424'; E='  +5.06%  ' },
  @{ Row=26; B='WrappedeETH'; C='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D='2.738.07'; E='  +8.39%  ' },
  @{ Row=27; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='''0.167'; E='  +3.64%  ' },
  @{ Row=28; B='Binance-PegBSC-USD'; C='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D='''0.990'; E='  -0.69%  ' },
  @{ Row=29; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0₃0863'; E='  +10.74%  ' },
  @{ Row=30; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='''7.57'; E='  +4.31%  ' },
  @{ Row=31; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='''0.997'; E='  -0.18%  ' },
  @{ Row=32; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='''19.52'; E='  +5.03%  ' },
  @{ Row=33; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='''156.18'; E='  +3.23%  ' },
  @{ Row=34; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='''1.57'; E='  +3.44%  ' },
  @{ Row=35; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='''5.74'; E='  +8.75%  ' },
  @{ Row=36; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='''4.04'; E='  +7.80%  ' },
  @{ Row=37; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='''1.21'; E='  +5.18%  ' },
  @{ Row=38; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='''310.29'; E='  +11.24%  ' },
  @{ Row=39; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='''1.49'; E='  +9.06%  ' },
  @{ Row=40; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='''0.851'; E='  +3.57%  ' },
  @{ Row=41; B='SuiNetwork'; C='https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; D='''0.839'; E='  +29.22%  ' },
  @{ Row=42; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='''3.76'; E='  +7.11%  ' },
  @{ Row=43; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='''35.41'; E='  +4.21%  ' },
  @{ Row=44; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='''0.633'; E='  +6.10%  ' },
  @{ Row=45; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='''0.0575'; E='  +8.53%  ' },
  @{ Row=46; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='''0.100'; E='  -1.50%  ' },
  @{ Row=47; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='''20.06'; E='  +14.27%  ' },
  @{ Row=48; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='''0.996'; E='  +0.10%  ' },
  @{ Row=49; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='''4.92'; E='  +7.65%  ' },
  @{ Row=50; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='''0.0236'; E='  +3.97%  ' },
  @{ Row=51; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='2.053.22'; E='  +8.10%  ' }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value2 = $item.B
    $ws.Cells.Item($r, 3).Value2 = $item.C
    $ws.Cells.Item($r, 4).Value2 = $item.D
    $ws.Cells.Item($r, 5).Value2 = $item.E
}

Write-Output "done"
